# Insert two new rows into the "general" sheet (rows 6 and 7), describing
# the prior distributions for fluxes and thermodynamic quantities, pushing
# all the rows below down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert two rows at row 6 (shifts existing rows 6..12 down to 8..14)
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "Prior distribution for fluxes (uniform or normal)"
$ws.Range("B6").Value = "normal"

$ws.Range("A7").Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Range("B7").Value = "normal"

# Match formatting: column A label cells look like the other label rows
# above (A2:A5 - bold Arial, boxed, left/top aligned); column B value
# cells get a boxed, centered look (default font).
$ws.Range("A2").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B6:B7").Borders.LineStyle = 1
$ws.Range("B6:B7").HorizontalAlignment = -4108  # xlHAlignCenter
$ws.Range("B6:B7").VerticalAlignment = -4107    # xlVAlignBottom

$ws.Activate()
$ws.Range("A6:B7").Select()
